$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 168 (shifts existing rows 168..383 down to 169..384)
$ws.Rows.Item(168).Insert()

# Populate the new row 168 with its data (matches the template used by sibling rows,
# with its own Fecha (D) and Volumen (J) values)
$ws.Range("A168").Value = 3
$ws.Range("B168").Value = "Femacal de La Calera"
$ws.Range("C168").Value = "Coquimbo"
$ws.Range("D168").Value = 44902
$ws.Range("E168").Value = 5
$ws.Range("F168").Value = 100112039
$ws.Range("G168").Value = "Ciboulette"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 160
$ws.Range("K168").Value = 1500
$ws.Range("L168").Value = 1500
$ws.Range("M168").Value = 1500
$ws.Range("N168").Value = '$/docena de atados'
$ws.Range("O168").Value = "Provincia de Quillota"
$ws.Range("P168").Value = 500
$ws.Range("Q168").Value = 3
$ws.Range("R168").Value = "Hortaliza"
